# Updated symbol list on Sat Dec 24 01:44:41 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    # Force the cell to be stored/treated as text, matching the source
    # workbook's inlineStr cells (avoids numeric auto-coercion of values
    # like "245.52" or "1").
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" '245.52'
Set-TextValue $ws "G2" '1'
Set-TextValue $ws "D3" '22.13'
Set-TextValue $ws "G3" '1'
Set-TextValue $ws "D4" '5.346'
Set-TextValue $ws "G4" '1'
Set-TextValue $ws "D5" '0.05889'
Set-TextValue $ws "G5" '1'
Set-TextValue $ws "D6" '3.394'
Set-TextValue $ws "G6" '1'
Set-TextValue $ws "D7" '6.382'
Set-TextValue $ws "G7" '1'
Set-TextValue $ws "D8" '0.8107'
Set-TextValue $ws "G8" '1'
Set-TextValue $ws "D9" '0.9629'
Set-TextValue $ws "G9" '1'
Set-TextValue $ws "B10" 'One'
Set-TextValue $ws "C10" 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue $ws "D10" '0.01120'
Set-TextValue $ws "E10" '9OneONEBestin24h'
Set-TextValue $ws "G10" '1'
Set-TextValue $ws "B11" 'WazirX'
Set-TextValue $ws "C11" 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws "D11" '0.1431'
Set-TextValue $ws "E11" '10WazirXWRX'
Set-TextValue $ws "G11" '1'
Set-TextValue $ws "D12" '0.07430'
Set-TextValue $ws "G12" '1'
Set-TextValue $ws "B13" 'LiechtensteinCryptoassetsExchange'
Set-TextValue $ws "C13" 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws "D13" '0.03464'
Set-TextValue $ws "E13" '12LiechtensteinCryptoassetsExchangeLCX'
Set-TextValue $ws "G13" '1'
Set-TextValue $ws "B14" 'BitrueCoin'
Set-TextValue $ws "C14" 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws "D14" '0.03034'
Set-TextValue $ws "E14" '13BitrueCoinBTR'
Set-TextValue $ws "G14" '1'
Set-TextValue $ws "B15" 'MCDex'
Set-TextValue $ws "C15" 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue $ws "D15" '4.422'
Set-TextValue $ws "E15" '14MCDexMCB'
Set-TextValue $ws "G15" '1'
Set-TextValue $ws "B16" 'BitMartToken'
Set-TextValue $ws "C16" 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws "D16" '0.09396'
Set-TextValue $ws "E16" '15BitMartTokenBMX'
Set-TextValue $ws "G16" '1'
Set-TextValue $ws "B17" 'BitForexToken'
Set-TextValue $ws "C17" 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws "D17" '0.001588'
Set-TextValue $ws "E17" '16BitForexTokenBF'
Set-TextValue $ws "G17" '1'
Set-TextValue $ws "B18" 'CoinExToken'
Set-TextValue $ws "C18" 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue $ws "D18" '0.04819'
Set-TextValue $ws "E18" '17CoinExTokenCET'
Set-TextValue $ws "G18" '1'
Set-TextValue $ws "D19" '0.006230'
Set-TextValue $ws "G19" '1'
Set-TextValue $ws "D20" '0.004078'
Set-TextValue $ws "G20" '1'
Set-TextValue $ws "D21" '0.0009867'
Set-TextValue $ws "G21" '1'
Set-TextValue $ws "D22" '0.00009705'
Set-TextValue $ws "G22" '1'
Set-TextValue $ws "D23" '3.700'
Set-TextValue $ws "G23" '1'
Set-TextValue $ws "D24" '2.216'
Set-TextValue $ws "G24" '1'
Set-TextValue $ws "D25" '0.3268'
Set-TextValue $ws "G25" '1'
Set-TextValue $ws "G26" '1'
Set-TextValue $ws "D27" '0.0002462'
Set-TextValue $ws "G27" '1'
Set-TextValue $ws "G28" '1'
Set-TextValue $ws "G29" '1'
Set-TextValue $ws "G30" '1'
Set-TextValue $ws "G31" '1'
Set-TextValue $ws "G32" '1'
Set-TextValue $ws "G33" '1'
Set-TextValue $ws "G34" '1'
Set-TextValue $ws "G35" '1'
Set-TextValue $ws "G36" '1'
Set-TextValue $ws "G37" '1'
Set-TextValue $ws "G38" '1'
Set-TextValue $ws "G39" '1'
Set-TextValue $ws "D40" '0.03927'
Set-TextValue $ws "G40" '1'
Set-TextValue $ws "B41" 'KickToken'
Set-TextValue $ws "C41" 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue $ws "D41" '0.006624'
Set-TextValue $ws "E41" '40KickTokenKICK'
Set-TextValue $ws "G41" '1'
Set-TextValue $ws "B42" 'BKEXToken'
Set-TextValue $ws "C42" 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue $ws "D42" '0.1073'
Set-TextValue $ws "E42" '41BKEXTokenBKK'
Set-TextValue $ws "G42" '1'
Set-TextValue $ws "B43" 'CEJI'
Set-TextValue $ws "C43" 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue $ws "D43" '0.003002'
Set-TextValue $ws "E43" '42CEJICEJI'
Set-TextValue $ws "G43" '1'
Set-TextValue $ws "G44" '1'
Set-TextValue $ws "D45" '0.00005303'
Set-TextValue $ws "G45" '1'
Set-TextValue $ws "G46" '1'
Set-TextValue $ws "D47" '0.6904'
Set-TextValue $ws "G47" '1'
Set-TextValue $ws "D48" '0.05475'
Set-TextValue $ws "E48" '47BOLOBOLOWorstin24h'
Set-TextValue $ws "G48" '1'
Set-TextValue $ws "G49" '1'
Set-TextValue $ws "D50" '0.01011'
Set-TextValue $ws "G50" '1'
Set-TextValue $ws "G51" '1'
